# Restructure solution projects for POM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet.
$ws.Name = "Browser"

# Clear the old data (A1:A8) and write the new, reordered values
# (PHANTOMJS_HEADLESS removed, browser order rearranged).
$ws.Range("A1:A8").ClearContents()

$values = @("BrowserName", "EDGE", "CHROME", "FIREFOX", "IE", "CHROME_HEADLESS", "FIREFOX_HEADLESS")
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Match the saved selection/active cell from the diff.
$ws.Range("A7").Select()
